# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list is reordered from descending (2208..2201) to
# ascending (2201..2208), and the "Salario Basico" (column G) is updated
# for every worker/period row from 689455 to 908526. The "Valor Mora"
# (column F) stays attached to its period: period 2208 keeps 26650, all
# the other periods keep 36341 - only its row position moves because the
# period ordering flipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16-23
$ws.Range("E16").Value = "2201"
$ws.Range("E17").Value = "2202"
$ws.Range("E18").Value = "2203"
$ws.Range("E19").Value = "2204"
$ws.Range("E20").Value = "2205"
$ws.Range("E21").Value = "2206"
$ws.Range("E22").Value = "2207"
$ws.Range("E23").Value = "2208"

# Valor Mora (column F) follows the period: 2208 -> 26650, rest -> 36341
$ws.Range("F16").Value = 36341
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 36341
$ws.Range("F22").Value = 36341
$ws.Range("F23").Value = 26650

# Salario Basico (column G) updated uniformly for every row
$ws.Range("G16:G23").Value = 908526
